$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 3 new (empty) columns before the current column F
#    (F=Valor_producto, G=Part_profesional, H=Revisar)
#    so that F,G,H become new blank columns and the old F,G,H shift to I,J,K
$ws.Range("F1:H1").EntireColumn.Insert()

# 2) Insert 1 more new (empty) column before the column that now holds
#    "Part_profesional" (which after step 1 is column J), shifting
#    Part_profesional/Revisar from J/K to K/L
$ws.Range("J1").EntireColumn.Insert()

# 3) Set the new header labels
$ws.Range("F1").Value = "Porc_trans"
$ws.Range("G1").Value = "Cost_trans"
$ws.Range("H1").Value = "Porc_producto"
$ws.Range("J1").Value = "Valor_Neto"

# 4) Fill in the new data values for rows 2-4
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.26
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.14575
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.26
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
